# Insert a new data row at row 33 (pushing the existing rows 33..139 down
# to 34..140) and populate it with the new "Poroto granado" price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(33).Insert()

$ws.Range("A33").Value = 5
$ws.Range("B33").Value = "Macroferia Regional de Talca"
$ws.Range("C33").Value = "Maule"
$ws.Range("D33").Value = 44648
$ws.Range("E33").Value = 7
$ws.Range("F33").Value = 100112030
$ws.Range("G33").Value = "Poroto granado"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 300
$ws.Range("K33").Value = 20000
$ws.Range("L33").Value = 20000
$ws.Range("M33").Value = 20000
$ws.Range("N33").Value = "`$/saco 25 kilos"
$ws.Range("O33").Value = "Región del Maule"
$ws.Range("P33").Value = 800
$ws.Range("Q33").Value = 25
$ws.Range("R33").Value = "Hortaliza"
